# The commit renames the header of column A from "County" to "county"
# (lowercase). All other header/data cells are unaffected; they simply
# point at shared-string indices that shift because the old "County"
# string is dropped from the shared-string table and "county" is
# appended as a new entry - that bookkeeping happens automatically
# when we just change the cell's displayed text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "county"
